$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing hyperlinks do not auto-shift when rows are inserted, so drop them
# first and rebuild the full set (in the new row order) at the end.
$ws.Hyperlinks.Delete()

# Insert a new row at the top of the data (row 2), pushing existing rows down.
$ws.Range("A2:F2").Insert()

# Copy formatting (only, not values) from the row below (now row 3, old row 2)
# into the freshly inserted row, matching just the used columns A:F.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row's values
$ws.Cells.Item(2, 1).Value = 8
$ws.Cells.Item(2, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(2, 3).Value = "IE07"
$ws.Cells.Item(2, 4).Value = 270.25
$ws.Cells.Item(2, 5).Value = "25-09-2025"
$ws.Cells.Item(2, 6).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf"

# Rebuild all hyperlinks for F2:F9 in order, top (newest) to bottom (oldest).
$urls = @(
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
)
for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = 2 + $i
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $urls[$i]) | Out-Null
    # Adding a hyperlink auto-applies Excel's built-in blue/underline "Hyperlink"
    # cell style; the source data only carries the link via <hyperlinks> (no
    # distinct cell style), so restore the plain formatting used by the rest
    # of the row right away.
    $excel.CutCopyMode = $false
    $ws.Cells.Item($row, 5).Copy()
    $ws.Cells.Item($row, 6).PasteSpecial(-4122) # xlPasteFormats
    $excel.CutCopyMode = $false
}
